# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (the stock Office palette)
#   ppt/theme/theme2.xml  -> "Integral"     (the palette actually driving
#                                             every slide, via slideMaster1)
#
# The authored edit swaps the two themes' contents in place (filenames and
# every relationship stay put) so the deck that used to render with the
# "Integral" palette now renders with the plain "Office" palette (and vice
# versa for the otherwise-unused notes-master theme).
#
# The PowerPoint object model doesn't let automation rename/re-target a
# theme part, but it does let us rewrite the 12 theme colour slots of the
# theme that is actually in force for the presentation (exposed through
# any slide's ThemeColorScheme, since it's the one shared part every slide
# / layout / master resolves to). Driving those 12 slots to the "Office"
# palette values reproduces the part of the swap that is visible anywhere
# in the presentation.

function ToThemeRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the stock "Office" colour scheme (previously theme1.xml),
# applied in clrScheme slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officePalette = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officePalette.Length; $i++) {
    $tcs.Item($i).RGB = ToThemeRGB($officePalette[$i - 1])
}
